# Regenerate merged AHB files
# - rename the "_old"/"_new" suffixed headers to "_FV2210"/"_FV2304"
# - freeze the header row
# - wrap the data range in an Excel table (ListObject)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

# Columns A-J (1-10) carried the "_old" suffix -> becomes "_FV2210"
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($headers[$i])_FV2210"
}

# Columns L-U (12-21) carried the "_new" suffix -> becomes "_FV2304"
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($headers[$i])_FV2304"
}

# Freeze the top header row (row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the used range into a proper Excel table
$range = $ws.Range("A1:U79")
$lo = $wb.ListObjects.Add(1, $range, 0, 1)
$lo.Name = "Table1"
